# Updates to dictionaries and moving pattern files
#
# The "Authority" column (G) on the Definitions worksheet referenced the
# shared string "AASHTO D 145" for every data row. That citation is updated
# to the correct AASHTO standard reference, "AASHTO M 145-91", for all rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definitions")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq "AASHTO D 145") {
        $cell.Value = "AASHTO M 145-91"
    }
}

# Restore the active sheet / selection state to match the authored workbook.
$ws.Activate()
$ws.Range("E15").Select()
